# Auto-generated edit script: updates Leve price/profit figures across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7999.6665
$ws.Range("I32").Value = 7999.6665
$ws.Range("K32").Value = 7999.6665
$ws.Range("M32").Value = -7673.6665
$ws.Range("H40").Value = 2400.2
$ws.Range("I40").Value = 2001
$ws.Range("K40").Value = 2001
$ws.Range("M40").Value = -1826
$ws.Range("H43").Value = 7999.4
$ws.Range("I43").Value = 9749.25
$ws.Range("K43").Value = 9749.25
$ws.Range("M43").Value = -9680.25
$ws.Range("H64").Value = 2600
$ws.Range("I64").Value = 2000
$ws.Range("K64").Value = 2000
$ws.Range("M64").Value = -1752
$ws.Range("H67").Value = 2600
$ws.Range("I67").Value = 2000
$ws.Range("K67").Value = 2000
$ws.Range("M67").Value = -1142
$ws.Range("H74").Value = 4354.048
$ws.Range("I74").Value = 3299.6667
$ws.Range("K74").Value = 3299.6667
$ws.Range("M74").Value = -2363.6667
$ws.Range("H76").Value = 7227.4546
$ws.Range("I76").Value = 6000.5
$ws.Range("K76").Value = 6000.5
$ws.Range("M76").Value = -5685.5
$ws.Range("H77").Value = 4354.048
$ws.Range("I77").Value = 3299.6667
$ws.Range("K77").Value = 16498.3335
$ws.Range("M77").Value = -11818.3335
$ws.Range("H79").Value = 7227.4546
$ws.Range("I79").Value = 6000.5
$ws.Range("K79").Value = 6000.5
$ws.Range("M79").Value = -4908.5
$ws.Range("H80").Value = 535
$ws.Range("J80").Value = 605
$ws.Range("L80").Value = 1815
$ws.Range("N80").Value = -3811
$ws.Range("H83").Value = 535
$ws.Range("J83").Value = 605
$ws.Range("L83").Value = 5445
$ws.Range("N83").Value = -15429
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 25000
$ws.Range("J34").Value = 25000
$ws.Range("L34").Value = 25000
$ws.Range("N34").Value = -25542
$ws.Range("H63").Value = 4918.857
$ws.Range("I63").Value = 1108
$ws.Range("K63").Value = 1108
$ws.Range("M63").Value = -422
$ws.Range("H66").Value = 4918.857
$ws.Range("I66").Value = 1108
$ws.Range("K66").Value = 5540
$ws.Range("M66").Value = -2108
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41748
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -128736
$ws.Range("H96").Value = 29500
$ws.Range("J96").Value = 29500
$ws.Range("L96").Value = 29500
$ws.Range("N96").Value = -34992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null
$ws.Range("H35").Value = 19285.285
$ws.Range("J35").Value = 19285.285
$ws.Range("L35").Value = 19285.285
$ws.Range("N35").Value = -19905.285
$ws.Range("H86").Value = 1668357.5
$ws.Range("I86").Value = 1399.6666
$ws.Range("J86").Value = 3335315.2
$ws.Range("K86").Value = 1399.6666
$ws.Range("L86").Value = 3335315.2
$ws.Range("M86").Value = -276.6666
$ws.Range("N86").Value = -3337561.2
$ws.Range("H89").Value = 1668357.5
$ws.Range("I89").Value = 1399.6666
$ws.Range("J89").Value = 3335315.2
$ws.Range("K89").Value = 6998.333000000001
$ws.Range("L89").Value = 16676576
$ws.Range("M89").Value = -1382.333000000001
$ws.Range("N89").Value = -16687808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -207
$ws.Range("N29").Value = $null
$ws.Range("H70").Value = 40000
$ws.Range("J70").Value = 40000
$ws.Range("L70").Value = 40000
$ws.Range("N70").Value = -40630
$ws.Range("H73").Value = 40000
$ws.Range("J73").Value = 40000
$ws.Range("L73").Value = 40000
$ws.Range("N73").Value = -42184
$ws.Range("H132").Value = 3800.2593
$ws.Range("I132").Value = 2267.1667
$ws.Range("K132").Value = 6801.500100000001
$ws.Range("M132").Value = -4271.500100000001
$ws.Range("H134").Value = 1727.1875
$ws.Range("I134").Value = 1293.0834
$ws.Range("K134").Value = 3879.2502
$ws.Range("M134").Value = -1344.2502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1666.6666
$ws.Range("I32").Value = 1666.6666
$ws.Range("K32").Value = 4999.9998
$ws.Range("M32").Value = -4716.9998
$ws.Range("H38").Value = 112.63158
$ws.Range("I38").Value = 57.266666
$ws.Range("J38").Value = 320.25
$ws.Range("K38").Value = 171.799998
$ws.Range("L38").Value = 960.75
$ws.Range("M38").Value = 175.200002
$ws.Range("N38").Value = -1654.75
$ws.Range("H94").Value = 2700
$ws.Range("H106").Value = 12666.667
$ws.Range("J106").Value = 12666.667
$ws.Range("L106").Value = 38000.001
$ws.Range("N106").Value = -39892.001
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null
$ws.Range("H131").Value = 2401.8147
$ws.Range("J131").Value = 2546.5208
$ws.Range("L131").Value = 7639.562399999999
$ws.Range("N131").Value = -17719.5624
$ws.Range("H140").Value = 4124.5
$ws.Range("I140").Value = 3833.3333
$ws.Range("K140").Value = 11499.9999
$ws.Range("M140").Value = -6319.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 3343004.2
$ws.Range("I35").Value = 3343004.2
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 3343004.2
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -3342706.2
$ws.Range("N35").Value = $null
$ws.Range("H39").Value = 55261
$ws.Range("J39").Value = 55261
$ws.Range("L39").Value = 55261
$ws.Range("N39").Value = -56325

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9549.125
$ws.Range("I16").Value = 9115.833000000001
$ws.Range("K16").Value = 9115.833000000001
$ws.Range("M16").Value = -8945.833000000001
$ws.Range("H40").Value = 1992.2307
$ws.Range("I40").Value = 1990.909
$ws.Range("J40").Value = 1999.5
$ws.Range("K40").Value = 1990.909
$ws.Range("L40").Value = 1999.5
$ws.Range("M40").Value = -1854.909
$ws.Range("N40").Value = -2271.5
$ws.Range("H68").Value = 4024.875
$ws.Range("I68").Value = 3699.8462
$ws.Range("J68").Value = 5433.3335
$ws.Range("K68").Value = 3699.8462
$ws.Range("L68").Value = 5433.3335
$ws.Range("M68").Value = -2950.8462
$ws.Range("N68").Value = -6931.3335
$ws.Range("H71").Value = 4024.875
$ws.Range("I71").Value = 3699.8462
$ws.Range("J71").Value = 5433.3335
$ws.Range("K71").Value = 18499.231
$ws.Range("L71").Value = 27166.6675
$ws.Range("M71").Value = -14755.231
$ws.Range("N71").Value = -34654.6675
$ws.Range("H122").Value = 3161.25
$ws.Range("I122").Value = 3166.6667
$ws.Range("J122").Value = 3145
$ws.Range("K122").Value = 9500.000100000001
$ws.Range("L122").Value = 9435
$ws.Range("M122").Value = -7050.000100000001
$ws.Range("N122").Value = -14335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 45005.5
$ws.Range("J20").Value = 45005.5
$ws.Range("L20").Value = 45005.5
$ws.Range("N20").Value = -45485.5
$ws.Range("H31").Value = 21947.6
$ws.Range("J31").Value = 21947.6
$ws.Range("L31").Value = 21947.6
$ws.Range("N31").Value = -22643.6
$ws.Range("H122").Value = 1223.1052
$ws.Range("I122").Value = 1223.1052
$ws.Range("K122").Value = 3669.3156
$ws.Range("M122").Value = -1219.3156
